$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.040.12'
$ws.Range('E2').Value = '  -1.77%  '
$ws.Range('D3').Value = '1.666.61'
$ws.Range('E3').Value = '  -1.42%  '
$ws.Range('D5').Value = "'216.18"
$ws.Range('E5').Value = '  -1.60%  '
$ws.Range('D6').Value = "'0.5094"
$ws.Range('E6').Value = '  -0.56%  '
$ws.Range('E7').Value = '  -0.20%  '
$ws.Range('D8').Value = "'0.2659"
$ws.Range('E8').Value = '  -0.49%  '
$ws.Range('D9').Value = "'0.06393"
$ws.Range('E9').Value = '  +1.22%  '
$ws.Range('D10').Value = "'21.88"
$ws.Range('E10').Value = '  -0.85%  '
$ws.Range('D11').Value = "'0.07461"
$ws.Range('E11').Value = '  +1.17%  '
$ws.Range('D12').Value = '1.688.78'
$ws.Range('E12').Value = '  -0.19%  '
$ws.Range('D13').Value = "'4.510"
$ws.Range('D14').Value = "'0.5809"
$ws.Range('E14').Value = '  +0.30%  '
$ws.Range('D15').Value = "'0.000008524"
$ws.Range('E15').Value = '  -0.84%  '
$ws.Range('E16').Value = '  -2.14%  '
$ws.Range('D17').Value = '26.115.32'
$ws.Range('E17').Value = '  -1.71%  '
$ws.Range('D18').Value = "'4.918"
$ws.Range('E18').Value = '  -1.39%  '
$ws.Range('D19').Value = "'1.005"
$ws.Range('E19').Value = '  -0.17%  '
$ws.Range('D20').Value = "'10.77"
$ws.Range('E20').Value = '  -1.46%  '
$ws.Range('D21').Value = "'189.94"
$ws.Range('E21').Value = '  +1.62%  '
$ws.Range('D22').Value = "'6.183"
$ws.Range('E22').Value = '  -1.39%  '
$ws.Range('E23').Value = '  -0.10%  '
$ws.Range('D24').Value = "'144.82"
$ws.Range('E24').Value = '  +0.03%  '
$ws.Range('D25').Value = "'7.599"
$ws.Range('E25').Value = '  +0.80%  '
$ws.Range('D26').Value = "'0.1204"
$ws.Range('E26').Value = '  +2.23%  '
$ws.Range('D27').Value = "'15.65"
$ws.Range('E27').Value = '  -0.78%  '
$ws.Range('D28').Value = "'0.06559"
$ws.Range('E28').Value = '  +12.33%  '
$ws.Range('E29').Value = '  -1.55%  '
$ws.Range('D30').Value = "'1.313"
$ws.Range('E30').Value = '  -2.14%  '
$ws.Range('D31').Value = "'3.553"
$ws.Range('E31').Value = '  +0.64%  '
$ws.Range('E32').Value = '  -0.45%  '
$ws.Range('D33').Value = "'1.658"
$ws.Range('E33').Value = '  +0.20%  '
$ws.Range('E34').Value = '  -0.48%  '
$ws.Range('E35').Value = '  +2.56%  '
$ws.Range('D36').Value = "'2.371"
$ws.Range('E36').Value = '  +0.31%  '
$ws.Range('D37').Value = "'2.687"
$ws.Range('E37').Value = '  +0.16%  '
$ws.Range('D38').Value = "'6.360"
$ws.Range('E38').Value = '  +8.06%  '
$ws.Range('D39').Value = '1.091.87'
$ws.Range('E39').Value = '  -0.77%  '
$ws.Range('D40').Value = "'0.01593"
$ws.Range('E40').Value = '  -1.66%  '
$ws.Range('D41').Value = "'0.8681"
$ws.Range('E41').Value = '  +0.23%  '
$ws.Range('E42').Value = '  +0.38%  '
$ws.Range('D43').Value = "'101.20"
$ws.Range('E43').Value = '  +1.42%  '
$ws.Range('D44').Value = '1.813.63'
$ws.Range('E44').Value = '  -1.88%  '
$ws.Range('D45').Value = "'0.00000000114"
$ws.Range('E45').Value = '  +0.25%  '
$ws.Range('D46').Value = "'56.29"
$ws.Range('E46').Value = '  -0.23%  '
$ws.Range('D47').Value = "'1.010"
$ws.Range('E47').Value = '  +0.23%  '
$ws.Range('D48').Value = "'8.081"
$ws.Range('E48').Value = '  +0.31%  '
$ws.Range('D49').Value = "'0.05227"
$ws.Range('E49').Value = '  -0.19%  '
$ws.Range('D50').Value = "'0.4289"
$ws.Range('E50').Value = '  -0.75%  '
$ws.Range('D51').Value = "'5.999"
$ws.Range('E51').Value = '  +2.90%  '
